$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) from 45204 to 45205 for all data rows (2-80)
for ($r = 2; $r -le 80; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45204) {
        $cell.Value2 = 45205
    }
}
